$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 153. This shifts the existing rows
# 153..253 down to 154..254 (carrying their data, including the date
# column values, with them) and keeps the sheet's overall structure/
# formatting intact.
$ws.Rows("153").Insert()

# Populate the newly inserted row 153 with the new weekly price record.
$ws.Range("A153").Value2 = 7
$ws.Range("B153").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C153").Value = "Ñuble"

$ws.Range("D153").Value2 = 44603
$ws.Range("D153").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("E153").Value2 = 16
$ws.Range("F153").Value2 = 100114013
$ws.Range("G153").Value = "Zanahoria"
$ws.Range("H153").Value = "Sin especificar"
$ws.Range("I153").Value = "Primera"
$ws.Range("J153").Value2 = 100
$ws.Range("K153").Value2 = 6000
$ws.Range("L153").Value2 = 6500
$ws.Range("M153").Value2 = 6250
$ws.Range("N153").Value = "$/saco 20 kilos"
$ws.Range("O153").Value = "Provincia de Diguillín"
$ws.Range("P153").Value2 = 312
$ws.Range("Q153").Value2 = 20
$ws.Range("R153").Value = "Hortaliza"
